$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting the latest cryptos data pull.
# Cells whose new text still parses as a plain number are pinned to
# Text format first so Excel keeps them as strings (matching the
# original inline-string cell type) instead of silently coercing
# them to numbers (which would also eat significant trailing zeros).
$ws.Range("D2").Value = "97.972.71"
$ws.Range("E2").Value = "  +3.16%  "
$ws.Range("D3").Value = "3.599.17"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.87"
$ws.Range("E5").Value = "  +3.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "657.51"
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.420"
$ws.Range("E8").Value = "  +6.44%  "
$ws.Range("E9").Value = "  +7.09%  "
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("D11").Value = "3.594.59"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.23"
$ws.Range("E12").Value = "  +5.68%  "
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "97.721.02"
$ws.Range("E15").Value = "  +3.06%  "
$ws.Range("D16").Value = "4.267.35"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("E17").Value = "  +3.67%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.63"
$ws.Range("E18").Value = "  +9.74%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.601.05"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.71"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.97"
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.519"
$ws.Range("E22").Value = "  +10.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.50"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "515.13"
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("E25").Value = "  +7.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.84"
$ws.Range("E26").Value = "  +3.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.06"
$ws.Range("E27").Value = "  +6.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.96"
$ws.Range("E28").Value = "  +5.54%  "
$ws.Range("D29").Value = "3.792.23"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.160"
$ws.Range("E30").Value = "  +15.14%  "
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.77"
$ws.Range("E32").Value = "  +4.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("E34").Value = "  +5.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.01"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.68"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("E37").Value = "  +7.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "621.79"
$ws.Range("E38").Value = "  +9.78%  "
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("E40").Value = "  +4.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.97"
$ws.Range("E41").Value = "  +14.52%  "
$ws.Range("E42").Value = "  +3.61%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.922"
$ws.Range("E44").Value = "  +3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.96"
$ws.Range("E45").Value = "  +7.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0442"
$ws.Range("E46").Value = "  +8.55%  "
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.67"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("B49").Value = "MantraDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.59"
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.50"
$ws.Range("E50").Value = "  +7.09%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.06"
$ws.Range("E51").Value = "  -0.58%  "
